# Applies the RWS "omreisroutes incl uitval onderliggend net" update to the
# "input" sheet of the RA2CE user_input_tests workbook:
#  - refreshes the shapefile ids (col A) and hazard_pickle file names (col P)
#    for the existing rows,
#  - appends 4 new analysis rows (13-16) with the same settings as the other
#    rows, and
#  - leaves the selection where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# --- Refresh existing rows 2-12: shapefile id (A) + hazard_pickle name (P) ---
$ws.Range("A2").Value = 13944
$ws.Range("P2").Value = "including_underlying/road_gdf_sel_incl_underl13944.p"
$ws.Range("A3").Value = 13946
$ws.Range("P3").Value = "including_underlying/road_gdf_sel_incl_underl13946.p"
$ws.Range("A4").Value = 70012
$ws.Range("P4").Value = "including_underlying/road_gdf_sel_incl_underl70012.p"
$ws.Range("A5").Value = 13814
$ws.Range("P5").Value = "including_underlying/road_gdf_sel_incl_underl13814.p"
$ws.Range("A6").Value = 13813
$ws.Range("P6").Value = "including_underlying/road_gdf_sel_incl_underl13813.p"
$ws.Range("A7").Value = 13943
$ws.Range("P7").Value = "including_underlying/road_gdf_sel_incl_underl13943.p"
$ws.Range("A8").Value = 70009
$ws.Range("P8").Value = "including_underlying/road_gdf_sel_incl_underl70009.p"
$ws.Range("A9").Value = 13165
$ws.Range("P9").Value = "including_underlying/road_gdf_sel_incl_underl13165.p"
$ws.Range("A10").Value = 14013
$ws.Range("P10").Value = "including_underlying/road_gdf_sel_incl_underl14013.p"
$ws.Range("A11").Value = 13173
$ws.Range("P11").Value = "including_underlying/road_gdf_sel_incl_underl13173.p"
$ws.Range("A12").Value = 19559
$ws.Range("P12").Value = "including_underlying/road_gdf_sel_incl_underl19559.p"

# --- Append new rows 13-16 (same settings as the rest of the table) ---
$ws.Range("A13:Y13").Borders.LineStyle = 1
$ws.Range("A13").Value = 13945
$ws.Range("B13").Value = $ws.Range("B12").Value()
$ws.Range("C13").Value = $ws.Range("C12").Value()
$ws.Range("D13").Value = $ws.Range("D12").Value()
$ws.Range("F13").Value = $ws.Range("F12").Value()
$ws.Range("L13").Value = $ws.Range("L12").Value()
$ws.Range("M13").Value = $ws.Range("M12").Value()
$ws.Range("N13").Value = $ws.Range("N12").Value()
$ws.Range("O13").Value = $ws.Range("O12").Value()
$ws.Range("P13").Value = "including_underlying/road_gdf_sel_incl_underl13945.p"
$ws.Range("R13").Value = $ws.Range("R12").Value()
$ws.Range("S13").Value = $ws.Range("S12").Value()
$ws.Range("T13").Value = 0.1

$ws.Range("A14:Y14").Borders.LineStyle = 1
$ws.Range("A14").Value = 13937
$ws.Range("B14").Value = $ws.Range("B12").Value()
$ws.Range("C14").Value = $ws.Range("C12").Value()
$ws.Range("D14").Value = $ws.Range("D12").Value()
$ws.Range("F14").Value = $ws.Range("F12").Value()
$ws.Range("L14").Value = $ws.Range("L12").Value()
$ws.Range("M14").Value = $ws.Range("M12").Value()
$ws.Range("N14").Value = $ws.Range("N12").Value()
$ws.Range("O14").Value = $ws.Range("O12").Value()
$ws.Range("P14").Value = "including_underlying/road_gdf_sel_incl_underl13937.p"
$ws.Range("R14").Value = $ws.Range("R12").Value()
$ws.Range("S14").Value = $ws.Range("S12").Value()
$ws.Range("T14").Value = 0.1

$ws.Range("A15:Y15").Borders.LineStyle = 1
$ws.Range("A15").Value = 13812
$ws.Range("B15").Value = $ws.Range("B12").Value()
$ws.Range("C15").Value = $ws.Range("C12").Value()
$ws.Range("D15").Value = $ws.Range("D12").Value()
$ws.Range("F15").Value = $ws.Range("F12").Value()
$ws.Range("L15").Value = $ws.Range("L12").Value()
$ws.Range("M15").Value = $ws.Range("M12").Value()
$ws.Range("N15").Value = $ws.Range("N12").Value()
$ws.Range("O15").Value = $ws.Range("O12").Value()
$ws.Range("P15").Value = "including_underlying/road_gdf_sel_incl_underl13812.p"
$ws.Range("R15").Value = $ws.Range("R12").Value()
$ws.Range("S15").Value = $ws.Range("S12").Value()
$ws.Range("T15").Value = 0.1

$ws.Range("A16:Y16").Borders.LineStyle = 1
$ws.Range("A16").Value = 19558
$ws.Range("B16").Value = $ws.Range("B12").Value()
$ws.Range("C16").Value = $ws.Range("C12").Value()
$ws.Range("D16").Value = $ws.Range("D12").Value()
$ws.Range("F16").Value = $ws.Range("F12").Value()
$ws.Range("L16").Value = $ws.Range("L12").Value()
$ws.Range("M16").Value = $ws.Range("M12").Value()
$ws.Range("N16").Value = $ws.Range("N12").Value()
$ws.Range("O16").Value = $ws.Range("O12").Value()
$ws.Range("P16").Value = "including_underlying/road_gdf_sel_incl_underl19558.p"
$ws.Range("R16").Value = $ws.Range("R12").Value()
$ws.Range("S16").Value = $ws.Range("S12").Value()
$ws.Range("T16").Value = 0.1

# --- Leave the selection where the author left it ---
$ws.Range("P28").Select()
